$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 2418462.5
$ws.Range("I76").Value = 3090.25
$ws.Range("J76").Value = 5053414
$ws.Range("K76").Value = 3090.25
$ws.Range("L76").Value = 5053414
$ws.Range("M76").Value = -2775.25
$ws.Range("N76").Value = -5054044

$ws.Range("H79").Value = 2418462.5
$ws.Range("I79").Value = 3090.25
$ws.Range("J79").Value = 5053414
$ws.Range("K79").Value = 3090.25
$ws.Range("L79").Value = 5053414
$ws.Range("M79").Value = -1998.25
$ws.Range("N79").Value = -5055598

$ws.Range("H129").Value = 859.43475
$ws.Range("J129").Value = 859.43475
$ws.Range("L129").Value = 2578.30425
$ws.Range("N129").Value = -12578.30425

$ws.Range("H137").Value = 46834.773
$ws.Range("I137").Value = 1064.2667
$ws.Range("K137").Value = 3192.800099999999
$ws.Range("M137").Value = -642.8000999999995

$ws.Range("H138").Value = 1727.0927
$ws.Range("I138").Value = 662.5238000000001
$ws.Range("J138").Value = 2404.5454
$ws.Range("K138").Value = 1987.5714
$ws.Range("L138").Value = 7213.6362
$ws.Range("M138").Value = 3152.4286
$ws.Range("N138").Value = -17493.6362

$ws.Range("H141").Value = 3200.7144
$ws.Range("I141").Value = 2439
$ws.Range("K141").Value = 7317
$ws.Range("M141").Value = -2137

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2332
$ws.Range("I2").Value = 1767
$ws.Range("K2").Value = 1767
$ws.Range("M2").Value = -1654

$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H32").Value = 18121.604
$ws.Range("I32").Value = 18720.797
$ws.Range("K32").Value = 18720.797
$ws.Range("M32").Value = -18433.797

$ws.Range("H63").Value = 4466428.5
$ws.Range("I63").Value = 2500
$ws.Range("K63").Value = 2500
$ws.Range("M63").Value = -1814

$ws.Range("H66").Value = 4466428.5
$ws.Range("I66").Value = 2500
$ws.Range("K66").Value = 12500
$ws.Range("M66").Value = -9068

$ws.Range("H74").Value = 43481132
$ws.Range("I74").Value = 55558424
$ws.Range("K74").Value = 55558424
$ws.Range("M74").Value = -55557550

$ws.Range("H77").Value = 43481132
$ws.Range("I77").Value = 55558424
$ws.Range("K77").Value = 277792120
$ws.Range("M77").Value = -277787752

$ws.Range("H116").Value = 2332
$ws.Range("I116").Value = 1767
$ws.Range("K116").Value = 1767
$ws.Range("M116").Value = 527

$ws.Range("H122").Value = 1711.4348
$ws.Range("I122").Value = 1707.45
$ws.Range("K122").Value = 5122.35
$ws.Range("M122").Value = -2672.35

$ws.Range("H132").Value = 19282.45
$ws.Range("I132").Value = 1779.2858
$ws.Range("K132").Value = 5337.857400000001
$ws.Range("M132").Value = -2807.857400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2332
$ws.Range("I3").Value = 1767
$ws.Range("K3").Value = 1767
$ws.Range("M3").Value = -1653

$ws.Range("H94").Value = 1807.78
$ws.Range("I94").Value = 918.0278
$ws.Range("K94").Value = 918.0278
$ws.Range("M94").Value = -467.0278

$ws.Range("H99").Value = 1584
$ws.Range("I99").Value = 1064.1666
$ws.Range("J99").Value = 2363.75
$ws.Range("K99").Value = 1064.1666
$ws.Range("L99").Value = 2363.75
$ws.Range("M99").Value = 433.8334
$ws.Range("N99").Value = -5359.75

$ws.Range("H105").Value = 3335647.2
$ws.Range("I105").Value = 2001.2858
$ws.Range("J105").Value = 6252587.5
$ws.Range("K105").Value = 2001.2858
$ws.Range("L105").Value = 6252587.5
$ws.Range("M105").Value = -254.2858000000001
$ws.Range("N105").Value = -6256081.5

$ws.Range("H107").Value = 850
$ws.Range("I107").Value = 700
$ws.Range("K107").Value = 700
$ws.Range("M107").Value = 1220

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 31.545454
$ws.Range("J7").Value = 29
$ws.Range("L7").Value = 29
$ws.Range("N7").Value = -255

$ws.Range("H16").Value = 1015.3333
$ws.Range("I16").Value = 1003.75
$ws.Range("J16").Value = 1028.5714
$ws.Range("K16").Value = 1003.75
$ws.Range("L16").Value = 1028.5714
$ws.Range("M16").Value = -716.75
$ws.Range("N16").Value = -1602.5714

$ws.Range("H31").Value = 16373.272
$ws.Range("I31").Value = 22794.2
$ws.Range("J31").Value = 2614.1428
$ws.Range("K31").Value = 22794.2
$ws.Range("L31").Value = 2614.1428
$ws.Range("M31").Value = -22499.2
$ws.Range("N31").Value = -3204.1428

$ws.Range("H34").Value = 16373.272
$ws.Range("I34").Value = 22794.2
$ws.Range("J34").Value = 2614.1428
$ws.Range("K34").Value = 22794.2
$ws.Range("L34").Value = 2614.1428
$ws.Range("M34").Value = -22592.2
$ws.Range("N34").Value = -3018.1428

$ws.Range("H94").Value = 7611.857
$ws.Range("I94").Value = 2313.1428
$ws.Range("J94").Value = 12910.571
$ws.Range("K94").Value = 2313.1428
$ws.Range("L94").Value = 12910.571
$ws.Range("M94").Value = -1862.1428
$ws.Range("N94").Value = -13812.571

$ws.Range("H113").Value = 1015.3333
$ws.Range("I113").Value = 1003.75
$ws.Range("J113").Value = 1028.5714
$ws.Range("K113").Value = 1003.75
$ws.Range("L113").Value = 1028.5714
$ws.Range("M113").Value = 1166.25
$ws.Range("N113").Value = -5368.5714

$ws.Range("H134").Value = 1044.7273
$ws.Range("I134").Value = 945.6
$ws.Range("K134").Value = 2836.8
$ws.Range("M134").Value = -301.8000000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 274.25
$ws.Range("I14").Value = 274.25
$ws.Range("K14").Value = 822.75
$ws.Range("M14").Value = -649.75

$ws.Range("H68").Value = 500
$ws.Range("J68").Value = 500
$ws.Range("L68").Value = 1500
$ws.Range("N68").Value = -3122

$ws.Range("H71").Value = 500
$ws.Range("J71").Value = 500
$ws.Range("L71").Value = 4500
$ws.Range("N71").Value = -12612

$ws.Range("H107").Value = 14507.143
$ws.Range("I107").Value = 50100
$ws.Range("J107").Value = 270
$ws.Range("K107").Value = 150300
$ws.Range("L107").Value = 810
$ws.Range("M107").Value = -148380
$ws.Range("N107").Value = -4650

$ws.Range("H131").Value = 660.33
$ws.Range("I131").Value = 566
$ws.Range("J131").Value = 665.29474
$ws.Range("K131").Value = 1698
$ws.Range("L131").Value = 1995.88422
$ws.Range("M131").Value = 3342
$ws.Range("N131").Value = -12075.88422

$ws.Range("H138").Value = 131952.39
$ws.Range("J138").Value = 251498.75
$ws.Range("L138").Value = 754496.25
$ws.Range("N138").Value = -764776.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 981.4643
$ws.Range("I102").Value = 976.9091
$ws.Range("K102").Value = 976.9091
$ws.Range("M102").Value = 645.0909

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4754.154
$ws.Range("I7").Value = 4410.4
$ws.Range("J7").Value = 5900
$ws.Range("K7").Value = 4410.4
$ws.Range("L7").Value = 5900
$ws.Range("M7").Value = -4298.4
$ws.Range("N7").Value = -6124

$ws.Range("H61").Value = 5624.55
$ws.Range("I61").Value = 2149.2
$ws.Range("J61").Value = 9099.9
$ws.Range("K61").Value = 2149.2
$ws.Range("L61").Value = 9099.9
$ws.Range("M61").Value = -1947.2
$ws.Range("N61").Value = -9503.9

$ws.Range("H68").Value = 4091.7896
$ws.Range("I68").Value = 2294.8
$ws.Range("J68").Value = 6088.4443
$ws.Range("K68").Value = 2294.8
$ws.Range("L68").Value = 6088.4443
$ws.Range("M68").Value = -1545.8
$ws.Range("N68").Value = -7586.4443

$ws.Range("H71").Value = 4091.7896
$ws.Range("I71").Value = 2294.8
$ws.Range("J71").Value = 6088.4443
$ws.Range("K71").Value = 11474
$ws.Range("L71").Value = 30442.2215
$ws.Range("M71").Value = -7730
$ws.Range("N71").Value = -37930.2215

$ws.Range("H100").Value = 3036.2273
$ws.Range("I100").Value = 2135.6428
$ws.Range("J100").Value = 4612.25
$ws.Range("K100").Value = 2135.6428
$ws.Range("L100").Value = 4612.25
$ws.Range("M100").Value = -1594.6428
$ws.Range("N100").Value = -5694.25

$ws.Range("H113").Value = 5624.55
$ws.Range("I113").Value = 2149.2
$ws.Range("J113").Value = 9099.9
$ws.Range("K113").Value = 2149.2
$ws.Range("L113").Value = 9099.9
$ws.Range("M113").Value = 20.80000000000018
$ws.Range("N113").Value = -13439.9

$ws.Range("H126").Value = 4754.154
$ws.Range("I126").Value = 4410.4
$ws.Range("J126").Value = 5900
$ws.Range("K126").Value = 13231.2
$ws.Range("L126").Value = 17700
$ws.Range("M126").Value = -10761.2
$ws.Range("N126").Value = -22640

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1766.1428
$ws.Range("I81").Value = 1766.1428
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 3532.2856
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -2471.2856
$ws.Range("N81").ClearContents()

$ws.Range("H82").Value = 20824.5
$ws.Range("J82").Value = 20824.5
$ws.Range("L82").Value = 20824.5
$ws.Range("N82").Value = -21590.5

$ws.Range("H84").Value = 1766.1428
$ws.Range("I84").Value = 1766.1428
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 17661.428
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -12357.428
$ws.Range("N84").ClearContents()

$ws.Range("H85").Value = 20824.5
$ws.Range("J85").Value = 20824.5
$ws.Range("L85").Value = 20824.5
$ws.Range("N85").Value = -23476.5

$ws.Range("H126").Value = 949.94446
$ws.Range("I126").Value = 825
$ws.Range("J126").Value = 1049.9
$ws.Range("K126").Value = 2475
$ws.Range("L126").Value = 3149.7
$ws.Range("M126").Value = -5
$ws.Range("N126").Value = -8089.700000000001
